# Insert a new data row before row 100 (pushing existing rows 100-206 down
# to 101-207) and populate the new row 100 with a fresh price observation.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("100:100").Insert()

$ws.Range("A100").Value2 = 10
$ws.Range("B100").Value2 = 'Vega Modelo de Temuco'
$ws.Range("C100").Value2 = 'La Araucanía'
$ws.Range("D100").Value2 = 44629
$ws.Range("E100").Value2 = 9
$ws.Range("F100").Value2 = 'Fruta'
$ws.Range("G100").Value2 = 100102
$ws.Range("H100").Value2 = 'Cítricos'
$ws.Range("I100").Value2 = 100102006
$ws.Range("J100").Value2 = 'Pomelo'
$ws.Range("K100").Value2 = 'Start Ruby'
$ws.Range("L100").Value2 = 'Primera'
$ws.Range("M100").Value2 = 30
$ws.Range("N100").Value2 = 15000
$ws.Range("O100").Value2 = 15000
$ws.Range("P100").Value2 = 15000
$ws.Range("Q100").Value2 = '$/bandeja 15 kilos granel'
$ws.Range("R100").Value2 = 'Región de O''Higgins'
$ws.Range("S100").Value2 = 1000
$ws.Range("T100").Value2 = 15
